# Actualización automática 2025-06-23 14:00:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M5").Value = 12.1
$ws1.Range("M22").Value = "4 de 20"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 12.1
$ws2.Range("F22").Value = 3409.33

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 3464.09
$ws3.Range("E16").Value = 26068.35
$ws3.Range("F16").Value = 0.1172977918519432
$ws3.Range("D19").Value = 3409.33
$ws3.Range("E19").Value = 46977.86762291769
$ws3.Range("F19").Value = 0.06766262385763899
